# Changes from the meeting
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update individual priority values
$ws.Range("B2").Value = "Very High"   # SRS
$ws.Range("B7").Value = "Medium"      # Country object
$ws.Range("B9").Value = "Medium"      # green wave

# Extend the dropdown validation list down to row 15 and add "Very High" option
$ws.Range("B2:B15").Validation.Delete()
$ws.Range("B2:B15").Validation.Add(3, 1, 1, '"Very High, High, Medium, Low,"')
$ws.Range("B2:B15").Validation.IgnoreBlank = $true
$ws.Range("B2:B15").Validation.InCellDropdown = $true
$ws.Range("B2:B15").Validation.ShowInput = $true
$ws.Range("B2:B15").Validation.ShowError = $true

# Update the selected cell in the sheet view
$ws.Range("C12").Select()
